$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns P (col_time) and Q (col_time_n) currently hold text representations
# of time-of-day values ("00:00:00" / "07:07:07"). Replace them with real
# numeric time values formatted as hh:mm:ss, matching the existing L/M/N/O
# date/datetime columns' treatment.

$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("P2:Q2").NumberFormat = "hh:mm:ss"

$ws.Range("P3").Value = 0.2966087962995516
$ws.Range("Q3").Value = 0.2966087962995516
$ws.Range("P3:Q3").NumberFormat = "hh:mm:ss"

$ws.Range("P4").Value = 0.2966087962995516
$ws.Range("P4").NumberFormat = "hh:mm:ss"
